$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 94
$ws.Range("I5").Value = 109.25
$ws.Range("J5").Value = 53.333332
$ws.Range("K5").Value = 109.25
$ws.Range("L5").Value = 53.333332
$ws.Range("M5").Value = 5.75
$ws.Range("N5").Value = -283.333332
$ws.Range("H6").Value = 186.18182
$ws.Range("I6").Value = 199.75
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 599.25
$ws.Range("L6").Value = 450
$ws.Range("M6").Value = -487.25
$ws.Range("N6").Value = -674
$ws.Range("H8").Value = 29.5
$ws.Range("I8").Value = 29
$ws.Range("K8").Value = 87
$ws.Range("M8").Value = 52
$ws.Range("H132").Value = 1957.9286
$ws.Range("I132").Value = 1957.9286
$ws.Range("K132").Value = 5873.7858
$ws.Range("M132").Value = -3343.7858
$ws.Range("H135").Value = 1408.5
$ws.Range("I135").Value = 1088.8
$ws.Range("J135").Value = 2207.75
$ws.Range("K135").Value = 9799.199999999999
$ws.Range("L135").Value = 19869.75
$ws.Range("M135").Value = -7264.199999999999
$ws.Range("N135").Value = -24939.75
$ws.Range("H137").Value = 1340.0454
$ws.Range("I137").Value = 1283.3158
$ws.Range("K137").Value = 3849.9474
$ws.Range("M137").Value = -1299.9474
$ws.Range("H138").Value = 2067.5398
$ws.Range("J138").Value = 2185
$ws.Range("L138").Value = 6555
$ws.Range("N138").Value = -16835

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 132.71428
$ws.Range("I4").Value = 122.5
$ws.Range("J4").Value = 146.33333
$ws.Range("K4").Value = 122.5
$ws.Range("L4").Value = 146.33333
$ws.Range("M4").Value = -6.5
$ws.Range("N4").Value = -378.33333
$ws.Range("H5").Value = 97.5
$ws.Range("I5").Value = 98.333336
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 98.333336
$ws.Range("L5").Value = 95
$ws.Range("M5").Value = 13.666664
$ws.Range("N5").Value = -319
$ws.Range("H102").Value = 2391
$ws.Range("I102").Value = 1883.75
$ws.Range("J102").Value = 2970.7144
$ws.Range("K102").Value = 1883.75
$ws.Range("L102").Value = 2970.7144
$ws.Range("M102").Value = -261.75
$ws.Range("N102").Value = -6214.7144

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 97.5
$ws.Range("I4").Value = 98.333336
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 98.333336
$ws.Range("L4").Value = 95
$ws.Range("M4").Value = 16.666664
$ws.Range("N4").Value = -325
$ws.Range("H86").Value = 2262.6155
$ws.Range("I86").Value = 2103.9
$ws.Range("K86").Value = 2103.9
$ws.Range("M86").Value = -980.9000000000001
$ws.Range("H89").Value = 2262.6155
$ws.Range("I89").Value = 2103.9
$ws.Range("K89").Value = 10519.5
$ws.Range("M89").Value = -4903.5
$ws.Range("H94").Value = 1759.1666
$ws.Range("I94").Value = 1022.4167
$ws.Range("K94").Value = 1022.4167
$ws.Range("M94").Value = -571.4167
$ws.Range("H99").Value = 1998.6666
$ws.Range("I99").Value = 1998.5
$ws.Range("K99").Value = 1998.5
$ws.Range("M99").Value = -500.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 942.7778
$ws.Range("I15").Value = 574.1667
$ws.Range("J15").Value = 1680
$ws.Range("K15").Value = 1722.5001
$ws.Range("L15").Value = 5040
$ws.Range("M15").Value = -1582.5001
$ws.Range("N15").Value = -5320
$ws.Range("H34").Value = 1256.8
$ws.Range("J34").Value = 1840.6
$ws.Range("L34").Value = 5521.799999999999
$ws.Range("N34").Value = -5689.799999999999
$ws.Range("H40").Value = 128.4
$ws.Range("I40").Value = 180.66667
$ws.Range("K40").Value = 722.66668
$ws.Range("M40").Value = -653.66668
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H128").Value = 2974693.5
$ws.Range("I128").Value = 2974693.5
$ws.Range("K128").Value = 8924080.5
$ws.Range("M128").Value = -8919100.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H70").Value = 6247.25
$ws.Range("I70").Value = 6663
$ws.Range("K70").Value = 6663
$ws.Range("M70").Value = -6393
$ws.Range("H73").Value = 6247.25
$ws.Range("I73").Value = 6663
$ws.Range("K73").Value = 6663
$ws.Range("M73").Value = -5727
$ws.Range("H132").Value = 2080.077
$ws.Range("I132").Value = 1731.375
$ws.Range("K132").Value = 5194.125
$ws.Range("M132").Value = -2664.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 982.9167
$ws.Range("I22").Value = 779.6
$ws.Range("J22").Value = 1128.1428
$ws.Range("K22").Value = 779.6
$ws.Range("L22").Value = 1128.1428
$ws.Range("M22").Value = -484.6
$ws.Range("N22").Value = -1718.1428
$ws.Range("H27").Value = 982.9167
$ws.Range("I27").Value = 779.6
$ws.Range("J27").Value = 1128.1428
$ws.Range("K27").Value = 779.6
$ws.Range("L27").Value = 1128.1428
$ws.Range("M27").Value = -672.6
$ws.Range("N27").Value = -1342.1428
$ws.Range("H55").Value = 1078.7778
$ws.Range("I55").Value = 957.4
$ws.Range("J55").Value = 1230.5
$ws.Range("K55").Value = 957.4
$ws.Range("L55").Value = 1230.5
$ws.Range("M55").Value = -784.4
$ws.Range("N55").Value = -1576.5
$ws.Range("H93").Value = 1078.0588
$ws.Range("I93").Value = 912.0833
$ws.Range("J93").Value = 1476.4
$ws.Range("K93").Value = 912.0833
$ws.Range("L93").Value = 1476.4
$ws.Range("M93").Value = 335.9167
$ws.Range("N93").Value = -3972.4
$ws.Range("H106").Value = 22129.75
$ws.Range("J106").Value = 22129.75
$ws.Range("L106").Value = 22129.75
$ws.Range("N106").Value = -24653.75
$ws.Range("H132").Value = 6600
$ws.Range("I132").Value = 6600
$ws.Range("K132").Value = 19800
$ws.Range("M132").Value = -17270

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 7614.2856
$ws.Range("I4").Value = 50000
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 50000
$ws.Range("L4").Value = 550
$ws.Range("M4").Value = -49887
$ws.Range("N4").Value = -776
$ws.Range("H96").Value = 5110.8
$ws.Range("J96").Value = 4531.5
$ws.Range("L96").Value = 4531.5
$ws.Range("N96").Value = -7277.5
$ws.Range("H100").Value = 2598
$ws.Range("I100").Value = 2337.7
$ws.Range("J100").Value = 3248.75
$ws.Range("K100").Value = 4675.4
$ws.Range("L100").Value = 6497.5
$ws.Range("M100").Value = -4134.4
$ws.Range("N100").Value = -7579.5
$ws.Range("H120").Value = 8000
$ws.Range("J120").Value = 8000
$ws.Range("L120").Value = 8000
$ws.Range("N120").Value = -17676
